$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text looks like a plain number need to be forced to
# text explicitly (NumberFormat "@"), otherwise Excel silently coerces
# them into a numeric value. The NumberFormat/Style dance avoids leaving
# a stray custom style behind once the text is in place.

$ws.Range("D2").Value = "27.489.69"
$ws.Range("E2").Value = "  -1.05%  "
$ws.Range("D3").Value = "1.834.51"
$ws.Range("E3").Value = "  -1.19%  "
$ws.Range("E4").Value = "  -2.77%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.15%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.005"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.57%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4307"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.14%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3710"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07282"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.83%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8684"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.25"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.65%  "
$ws.Range("D12").Value = "1.836.58"
$ws.Range("E12").Value = "  -1.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.700"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.374"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.07106"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.11%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.008"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008931"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.005"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.58%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "15.28"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.69%  "
$ws.Range("D21").Value = "27.492.61"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.182"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.86%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.12%  "
$ws.Range("D24").Value = "2.055.78"
$ws.Range("E24").Value = "  -1.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.006"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.160"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +7.60%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.298"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "117.78"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08880"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.38%  "
$ws.Range("E32").Value = "  -0.26%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7694"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.55%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.500"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.08%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.910"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.45%  "
$ws.Range("E36").Value = "  -2.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.125"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01964"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.93%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05292"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.00%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.190"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.16%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.881"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.84%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5097"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.91%  "
$ws.Range("B43").Value = "Algorand"
$ws.Range("C43").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1677"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.714"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.50%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.64"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.26%  "
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "106.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.87%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.06436"
$ws.Range("D48").Style = "Normal"
$ws.Range("E49").Value = "  -2.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.673"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.24%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.829"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.48%  "
